$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-9 (columns E, F, G)
$updates = @(
    @{Row=2;  E=1558; F=37316176; G=96447492},
    @{Row=3;  E=1567; F=38342720; G=112026199},
    @{Row=4;  E=1468; F=35654576; G=110866334},
    @{Row=5;  E=1564; F=38302544; G=88241881},
    @{Row=6;  E=1513; F=36836608; G=91482447},
    @{Row=7;  E=1521; F=36946080; G=98949647},
    @{Row=8;  E=1614; F=38319320; G=79650292},
    @{Row=9;  E=1536; F=37098904; G=102882794}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}

# Clear rows 10 and 11 entirely (A:G), leaving styles intact
$ws.Range("A10:G11").ClearContents()

# Set column F width (bestFit-like) to match new widest content
# (closest value reachable through Excel's quantized ColumnWidth property)
$ws.Columns.Item(6).ColumnWidth = 10.5

# Update selection to G9
$ws.Range("G9").Select()
